$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 339, shifting all existing rows (339..395) down to (340..396),
# then populate the new row 339 with this week's new price record.
$ws.Rows.Item(339).Insert()

$ws.Cells.Item(339, 1).Value = 10
$ws.Cells.Item(339, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(339, 3).Value = "La Araucanía"
$ws.Cells.Item(339, 4).Value = 44776
$ws.Cells.Item(339, 5).Value = 9
$ws.Cells.Item(339, 6).Value = 100112040
$ws.Cells.Item(339, 7).Value = "Cilantro"
$ws.Cells.Item(339, 8).Value = "Sin especificar"
$ws.Cells.Item(339, 9).Value = "Primera"
$ws.Cells.Item(339, 10).Value = 50
$ws.Cells.Item(339, 11).Value = 4000
$ws.Cells.Item(339, 12).Value = 4000
$ws.Cells.Item(339, 13).Value = 4000
$ws.Cells.Item(339, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(339, 15).Value = "Región Metropolitana"
$ws.Cells.Item(339, 16).Value = 2000
$ws.Cells.Item(339, 17).Value = 2
$ws.Cells.Item(339, 18).Value = "Hortaliza"
